$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) from the existing data row (row 3) down onto
# the two new rows being populated, matching the s="1" cell style used by
# the rest of the table.
$ws.Range("A3:C3").Copy()
$ws.Range("A4:C4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A5:C5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New row 4 values
$ws.Range("A4").Value = "Joberno "
$ws.Range("B4").Value = "j@gmail.com"
$ws.Range("C4").Value = 50

# New row 5 values
$ws.Range("A5").Value = "ahsah"
$ws.Range("B5").Value = "asdasd@gmail.com"
$ws.Range("C5").Value = 40
